$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.066.90"
$ws.Range("E2").Value = "  +1.31%  "

$ws.Range("D3").Value = "2.151.46"
$ws.Range("E3").Value = "  -0.57%  "

$ws.Range("E4").Value = "  -0.15%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "250.87"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +5.60%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.602"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.84%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "72.45"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.98%  "

$ws.Range("E8").Value = "  -0.05%  "

$ws.Range("E9").Value = "  -0.66%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.20"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.15%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0904"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.09%  "

$ws.Range("E12").Value = "  +0.57%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.68"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.11%  "

$ws.Range("D14").Value = "2.473.25"
$ws.Range("E14").Value = "  -0.61%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.05"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.53%  "

$ws.Range("D16").Value = "2.170.92"
$ws.Range("E16").Value = "  +1.35%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.759"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.18%  "

$ws.Range("D18").Value = "41.917.15"
$ws.Range("E18").Value = "  +1.30%  "

$ws.Range("E19").Value = "  -1.48%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "70.14"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.27%  "

$ws.Range("E21").Value = "  +0.32%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "224.93"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.52%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.49"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.31%  "

$ws.Range("E24").Value = "  +5.31%  "

$ws.Range("E25").Value = "  -0.11%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.37"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.91%  "

$ws.Range("E27").Value = "  +1.04%  "

$ws.Range("E28").Value = "  +2.94%  "

$ws.Range("E29").Value = "  -0.27%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "36.17"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +9.07%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "167.63"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.95%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "19.84"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.28%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0794"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.41%  "

$ws.Range("E34").Value = "  -3.92%  "

$ws.Range("E35").Value = "  -0.62%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.106"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.99%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.18"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.52%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0325"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +7.30%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "11.63"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.24%  "

$ws.Range("E40").Value = "  -3.05%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.194"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.41%  "

$ws.Range("B42").Value = "MultiversX"
$ws.Range("C42").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "58.15"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.11%  "

$ws.Range("B43").Value = "THORChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.08"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.00%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "100.66"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.86%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.459"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +14.87%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.16"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.58%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0957"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.83%  "

$ws.Range("E48").Value = "  +8.24%  "

$ws.Range("E49").Value = "  -0.06%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.11"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.22%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.63"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.42%  "
